# Updates Price (D) and Volume(1h) (E) columns for the crypto rows (2-51)
# on Sheet1, per the upstream coinranking.com refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Some of the new values are plain decimal numbers (e.g. "211.65").
    # Assigning those straight to .Value lets Excel auto-coerce them to
    # a Number, which silently drops significant trailing zeros (e.g.
    # "62.50" -> 62.5). Force text storage, write, then drop the
    # number-format override so the cell keeps its original (unstyled)
    # look, matching cells whose value was never at risk of coercion.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "28.521.05"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.563.60"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue $ws.Range("D5") "211.65"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("E7").Value = "  -0.05%  "
Set-TextValue $ws.Range("D8") "46.31"
$ws.Range("E8").Value = "  +5.29%  "
Set-TextValue $ws.Range("D9") "24.17"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").Value = "1.788.27"
$ws.Range("D14").Value = "1.569.47"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("D16").Value = "28.524.68"
$ws.Range("E16").Value = "  +0.14%  "
Set-TextValue $ws.Range("D17") "3.68"
$ws.Range("E17").Value = "  -3.08%  "
Set-TextValue $ws.Range("D18") "62.06"
$ws.Range("E18").Value = "  -3.01%  "
Set-TextValue $ws.Range("D19") "227.65"
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("E20").Value = "  -2.19%  "
Set-TextValue $ws.Range("D21") "7.31"
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("E22").Value = "  +0.00%  "
Set-TextValue $ws.Range("D23") "3.85"
$ws.Range("E23").Value = "  -6.84%  "
Set-TextValue $ws.Range("D24") "9.12"
$ws.Range("E24").Value = "  -3.09%  "
$ws.Range("E25").Value = "  +5.36%  "
Set-TextValue $ws.Range("D26") "150.35"
$ws.Range("E26").Value = "  -1.04%  "
Set-TextValue $ws.Range("D27") "14.96"
$ws.Range("E27").Value = "  -2.15%  "
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("E29").Value = "  -3.71%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("E33").Value = "  -1.68%  "
Set-TextValue $ws.Range("D34") "3.13"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("D35").Value = "1.396.47"
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("E37").Value = "  -3.53%  "
$ws.Range("E38").Value = "  +1.05%  "
Set-TextValue $ws.Range("D39") "2.57"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("E40").Value = "  -0.99%  "
Set-TextValue $ws.Range("D41") "0.535"
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  -3.88%  "
$ws.Range("E44").Value = "  +2.38%  "
Set-TextValue $ws.Range("D45") "5.51"
$ws.Range("E45").Value = "  -4.60%  "
$ws.Range("E46").Value = "  -0.32%  "
Set-TextValue $ws.Range("D47") "62.50"
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("D48").Value = "1.699.92"
$ws.Range("E48").Value = "  -1.89%  "
Set-TextValue $ws.Range("D49") "86.18"
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("E50").Value = "  -4.11%  "
$ws.Range("E51").Value = "  -0.87%  "
